$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "Jassem Al Memari"
$summary.Range("B4").Value = 2652.68
$summary.Range("B6").Value = 4992
$summary.Range("B7").Value = 87212
$summary.Range("B8").Value = -82220
$summary.Range("B9").Value = 0.06

# ---------------------------------------------------------------------
# Assets sheet
# Remove the "Mid-range Car" row (old row 3) and the old "Liquid Assets"
# row (old row 4), leaving the TOTAL ASSETS row right after row 2.
# ---------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")
$assets.Rows("3:4").Delete()

$assets.Range("A2").Value = "Liquid Assets"
$assets.Range("B2").Value = "Savings Account"
$assets.Range("C2").Value = 4992

$assets.Range("C3").Value = 4992

# ---------------------------------------------------------------------
# Liabilities sheet
# Remove the old "Personal Loan" (row 4) and "Credit Card Balance"
# (row 5) rows, since the Personal Loan / Credit Card data is being
# moved up into rows 2 and 3 (replacing the old Auto Loan rows).
# ---------------------------------------------------------------------
$liabilities = $wb.Worksheets.Item("Liabilities")
$liabilities.Rows("4:5").Delete()

$liabilities.Range("A2").Value = "Personal Loans"
$liabilities.Range("B2").Value = "Personal Loan"
$liabilities.Range("C2").Value = 61581
$liabilities.Range("D2").Value = 1283
$liabilities.Range("E2").Value = 4

$liabilities.Range("A3").Value = "Credit Cards"
$liabilities.Range("B3").Value = "Credit Card Balance"
$liabilities.Range("C3").Value = 25631
$liabilities.Range("D3").Value = 1282
$liabilities.Range("E3").Value = 1

$liabilities.Range("C4").Value = 87212

Write-Host "edits applied"
